# The "Employment impact" row (row 69) is removed from the Specification
# sheet. All subsequent rows (70-181) shift up by one to fill the gap,
# which also naturally updates the dependent row merges (column A/B
# section-header merge ranges) and the sheet's used-range dimension.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(69).Delete()
